$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-239 changes from serial date 45203 to 45205
$ws.Range("C2:C239").Value = 45205
